$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal:" value from EF-7 to EF-5
$ws.Range("B9").Value = "EF-5"
$ws.Range("C9").Value = "EF-5"

# Remove the "Requisitos:" rows (22 and 23)
$ws.Rows("22:23").Delete()
